$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C222").Value = 1069.578287
$ws.Range("C223").Value = 1080.6319803
$ws.Range("C224").Value = 1091.6835497
$ws.Range("C225").Value = 1102.7188191
$ws.Range("C226").Value = 1113.7230876
$ws.Range("C227").Value = 1124.6813456
$ws.Range("C228").Value = 1135.5785122
$ws.Range("C229").Value = 1146.3996904
$ws.Range("C230").Value = 1157.1304303
$ws.Range("C231").Value = 1167.7569945
$ws.Range("C232").Value = 1178.266614
$ws.Range("C233").Value = 1188.6477261
$ws.Range("C234").Value = 1198.8901854
$ws.Range("C235").Value = 1208.985438
$ws.Range("C236").Value = 1218.9266535
$ws.Range("C237").Value = 1228.7088082
$ws.Range("C238").Value = 1238.3287174
$ws.Range("C239").Value = 1247.785016
$ws.Range("C240").Value = 1257.078089
$ws.Range("C241").Value = 1266.2099555
$ws.Range("C242").Value = 1275.184114
$ws.Range("C243").Value = 1284.0053537
$ws.Range("C244").Value = 1292.6795427
$ws.Range("C245").Value = 1301.2134002
$ws.Range("C246").Value = 1309.614262
$ws.Range("C247").Value = 1317.8898471
$ws.Range("C248").Value = 1326.0480331
$ws.Range("C249").Value = 1334.0966439
$ws.Range("C250").Value = 1342.0432577
$ws.Range("C251").Value = 1349.8950344
$ws.Range("C252").Value = 1357.6585675
$ws.Range("C253").Value = 1365.3397597
$ws.Range("C254").Value = 1372.9437229
$ws.Range("C255").Value = 1380.4747028
$ws.Range("C256").Value = 1387.9360282
$ws.Range("C257").Value = 1395.3300834
$ws.Range("C258").Value = 1402.658306
$ws.Range("C259").Value = 1409.9212085
$ws.Range("C260").Value = 1417.1184255
$ws.Range("C261").Value = 1424.2487853
$ws.Range("C262").Value = 1431.3104091
$ws.Range("C263").Value = 1438.300834
$ws.Range("C264").Value = 1445.2171625
$ws.Range("C265").Value = 1452.0562341
$ws.Range("C266").Value = 1458.8148176
$ws.Range("C267").Value = 1465.4898201
$ws.Range("C268").Value = 1472.0785062
$ws.Range("C269").Value = 1478.5787216
$ws.Range("C270").Value = 1484.9891145
$ws.Range("C271").Value = 1491.3093441
$ws.Range("C272").Value = 1497.5402706
$ws.Range("C273").Value = 1503.6841168
$ws.Range("C274").Value = 1509.7445937
$ws.Range("C275").Value = 1515.7269841
$ws.Range("C276").Value = 1521.6381781
$ws.Range("C277").Value = 1527.4866555
$ws.Range("C278").Value = 1533.2824131
$ws.Range("C279").Value = 1539.0368331
$ws.Range("C280").Value = 1544.7624907
$ws.Range("C281").Value = 1550.4729005
$ws.Range("C282").Value = 1556.1821988
$ws.Range("C283").Value = 1561.9047625
$ws.Range("C284").Value = 1567.6547623
$ws.Range("C285").Value = 1573.4456546
$ws.Range("C286").Value = 1579.2896146
$ws.Range("C287").Value = 1585.1969217
$ws.Range("C288").Value = 1591.1753137
$ws.Range("C289").Value = 1597.2293352
$ws.Range("C290").Value = 1603.3597178
$ws.Range("C291").Value = 1609.5628387
$ws.Range("C292").Value = 1615.8303139
$ws.Range("C293").Value = 1622.1487902
$ws.Range("C294").Value = 1628.4999946
$ws.Range("C295").Value = 1634.861091
$ws.Range("C296").Value = 1641.2053705
$ws.Range("C297").Value = 1647.5032695
$ws.Range("C298").Value = 1653.723667
$ws.Range("C299").Value = 1659.8353736
$ws.Range("C300").Value = 1665.8086876
$ws.Range("C301").Value = 1671.6168765
$ws.Range("C302").Value = 1677.2374406
$ws.Range("C303").Value = 1682.6530403
$ws.Range("C304").Value = 1687.8520115
$ws.Range("C305").Value = 1692.828442
$ws.Range("C306").Value = 1697.5818396
$ws.Range("C307").Value = 1702.1164638
$ws.Range("C308").Value = 1706.44042
$ws.Range("C309").Value = 1710.5646298
$ws.Range("C310").Value = 1714.5017804
$ws.Range("C311").Value = 1718.2653417
$ws.Range("C312").Value = 1721.8687145
$ws.Range("C313").Value = 1725.324546
$ws.Range("C314").Value = 1728.6442289
$ws.Range("C315").Value = 1731.8375798
$ws.Range("C316").Value = 1734.91268
$ws.Range("C317").Value = 1737.8758573
$ws.Range("C318").Value = 1740.7317782
$ws.Range("C319").Value = 1743.4836283
$ws.Range("C320").Value = 1746.1333506
$ws.Range("C321").Value = 1748.6819236
$ws.Range("C322").Value = 1751.1296555
$ws.Range("C323").Value = 1753.4764807
$ws.Range("C324").Value = 1755.7222429
$ws.Range("C325").Value = 1757.8669521
$ws.Range("C326").Value = 1759.911007
$ws.Range("C327").Value = 1761.8553738
$ws.Range("C328").Value = 1763.701716
$ws.Range("C329").Value = 1765.4524707
$ws.Range("C330").Value = 1767.1108711
$ws.Range("C331").Value = 1768.6809151
$ws.Range("C332").Value = 1770.1672849
$ws.Range("C333").Value = 1771.5752228
$ws.Range("C334").Value = 1772.9103716
$ws.Range("C335").Value = 1774.1785908
$ws.Range("C336").Value = 1775.3857587
$ws.Range("C337").Value = 1776.5375739
$ws.Range("C338").Value = 1777.6393677
$ws.Range("C339").Value = 1778.6959382
$ws.Range("C340").Value = 1779.711418
$ws.Range("C341").Value = 1780.6891797
$ws.Range("C342").Value = 1781.6317868
$ws.Range("C343").Value = 1782.5409906
$ws.Range("C344").Value = 1783.4177705
$ws.Range("C345").Value = 1784.262415
$ws.Range("C346").Value = 1785.0746347
$ws.Range("C347").Value = 1785.8536987
$ws.Range("C348").Value = 1786.5985836
$ws.Range("C349").Value = 1787.3081252
$ws.Range("C350").Value = 1787.9811617
$ws.Range("C351").Value = 1788.616662
$ws.Range("C352").Value = 1789.2138295
